$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("A2").Value = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"
$overview.Range("G2").Value = "2016-08-17 03:02:04"

$hl = $overview.Hyperlinks.Item(1)
$hl.TextToDisplay = "e2e\ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("A2").Value = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"
$zhcn.Range("G2").Value = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.b31c831d896f5986db58330c22e5783791f12f69.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-17 03:01:56"
$zhcn.Range("I2").Value = ""
$zhcn.Range("J2").Value = ""
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

$zhcnHl = $zhcn.Hyperlinks.Item(1)
$zhcnHl.TextToDisplay = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"
$zhcn.Range("I2").Hyperlinks.Delete()

$zhcn.Columns.Item(9).ColumnWidth = 17.75
$zhcn.Columns.Item(10).ColumnWidth = 20.75

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("A2").Value = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"
$dede.Range("G2").Value = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.b31c831d896f5986db58330c22e5783791f12f69.de-de.xlf"
$dede.Range("H2").Value = "2016-08-17 03:02:04"
$dede.Range("I2").Value = ""
$dede.Range("J2").Value = ""
$dede.Range("K2").Value = "0001-01-01 00:00:00"

$dedeHl = $dede.Hyperlinks.Item(1)
$dedeHl.TextToDisplay = "ed303d24-600d-4b21-ad0d-f2834c8b0dcf.md"
$dede.Range("I2").Hyperlinks.Delete()

$dede.Columns.Item(9).ColumnWidth = 17.75
$dede.Columns.Item(10).ColumnWidth = 20.75
